$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new data rows right before the existing "Provincia del Elquí" block
# (old rows 26-42 shift down to 31-47), then fill the new rows 26-30 with the
# new "Provincia de Limarí" weekly data.
$ws.Rows("26:30").Insert()

# Shared/common values for every row in this block.
$mercadoId = 9
$mercado = "Vega Central Mapocho de Santiago"
$region = "Metropolitana"
$codreg = 13
$tipo = "Fruta"
$productoId = 100107
$producto = "Otros"
$categoriaId = 100107002
$categoria = "Chirimoya"
$variedad = "Cultivar IV Región"
$origen = "Provincia de Limarí"

function Set-Row($r, $fecha, $calidad, $volumen, $pmin, $pmax, $pprom, $unidad, $precioKg, $kgUnidad) {
    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $tipo
    $ws.Cells.Item($r, 7).Value = $productoId
    $ws.Cells.Item($r, 8).Value = $producto
    $ws.Cells.Item($r, 9).Value = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $calidad
    $ws.Cells.Item($r, 13).Value = $volumen
    $ws.Cells.Item($r, 14).Value = $pmin
    $ws.Cells.Item($r, 15).Value = $pmax
    $ws.Cells.Item($r, 16).Value = $pprom
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $precioKg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}

Set-Row 26 44463 "Especial"                220 24000 24000 24000 "`$/bandeja 8 kilos"              3000 8
Set-Row 27 44463 "Extra (doble especial)"   200 25600 25600 25600 "`$/bandeja 8 kilos"              3200 8
Set-Row 28 44463 "Primera"                  200 20000 20000 20000 "`$/bandeja 8 kilos"              2500 8
Set-Row 29 44463 "Segunda"                  250 16000 16000 16000 "`$/bandeja 8 kilos"              2000 8
Set-Row 30 44463 "Tercera"                  200 1500  1500  1500  "`$/kilo (en caja de 15 kilos)"   1500 1
